# Add the new YouTube-channel rows (28-59) to the "Sorted" worksheet,
# then leave the "Sorted" sheet active/selected (matching the final
# workbook state described by the commit).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Sorted"

$rows = @(
    @{Row=28; A='Web Development'; B='Net Ninjs'},
    @{Row=29; A='Python'; B='Corey Sohafer'},
    @{Row=30; A='Cloud'; B='CloudGuru'},
    @{Row=31; A='Cyber Security'; B='David Bombal'},
    @{Row=32; A='Quick Videos'; B='Fireship'},
    @{Row=33; A='AI/ML'; B='Andrei Karpathy'},
    @{Row=34; A='JavaScript'; B='Bro Code'},
    @{Row=35; A='CSS'; B='Kevin Powell'},
    @{Row=36; A='Algorithms'; B='Abdul Bari'},
    @{Row=37; A='Data Science'; B='StatQuest'},
    @{Row=38; A='Java'; B='Telusko'},
    @{Row=39; A='LeetCode'; B='NeetCode'},
    @{Row=40; A='Azure/AWS'; B='Andrew Brown'},
    @{Row=41; A='SQL'; B='Luke Barousse'},
    @{Row=42; A='Java'; B='Kunal Kushwaha'},
    @{Row=43; A='Blockchain'; B='Telusko'},
    @{Row=44; A='AI/ML'; B='Krish Naik'},
    @{Row=45; A='AR/VR'; B='FuseDVR'},
    @{Row=46; A='Chemistry'; B='Khan Academy'},
    @{Row=47; A='UI/UX'; B='GFXMentor'},
    @{Row=48; A='Devops'; B='Tech world with Nana'},
    @{Row=49; A='Cyber Security'; B='The Cyber Mentor'},
    @{Row=50; A='Ruby'; B='The Ruby Way'},
    @{Row=51; A='Scala'; B='Scala Love'},
    @{Row=52; A='JavaScript'; B='Traversy Media'},
    @{Row=53; A='Python'; B='Code With Harry'},
    @{Row=54; A='Kotlin'; B='Kotlin Programming'},
    @{Row=55; A='Flutter'; B='The Net Ninja'},
    @{Row=56; A='C'; B='FreeCodeCamp.org'},
    @{Row=57; A='C++'; B='The Cherno'},
    @{Row=58; A='SQL'; B='Programming With Mosh'},
    @{Row=59; A='Web Development'; B='Traversy Media'}
)

foreach ($r in $rows) {
    $ws1.Cells.Item($r.Row, 1).Value = $r.A
    $ws1.Cells.Item($r.Row, 2).Value = $r.B
}

# Make "Sorted" the active sheet/tab (was "Unsorted"), with the
# cursor resting near the newly-added data.
$ws1.Activate() | Out-Null
$ws1.Range("D48").Select() | Out-Null

Write-Host "Added" $rows.Count "rows to" $ws1.Name
